$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "নাম: Dr. Md. Shahjahan"
$ws.Range("A4").Value = "পদবী: ডিন"
$ws.Range("F5").Value = "বিভাগ :EEE অনুষদ"
$ws.Range("G26").Value = 5
$ws.Range("A32").Value = "কথায়:তের হাজার পাঁচশত টাকা মাত্র।"

$ws.Range("B5").Select()
